$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. The '_GoBack' bookmark (Word's "last edit location" marker) was
#    sitting at the start of the 'PeLib' heading. Move it so it sits
#    inside the 'Pattern generator.' bullet instead, right after 'Pat'
#    (i.e. splitting the run 'Pattern generator.' into 'Pat' + bookmark
#    + 'tern generator.').
# ------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$patRange = $d.Content
$patRange.Find.Execute("Pattern generator.", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0) | Out-Null
$splitPoint = $patRange.Start + 3   # after 'Pat'
$bookmarkRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

# ------------------------------------------------------------------
# 2. Remove the now-completed todo item "Custom file format for
#    specifying patterns." from the FindPattern list.
# ------------------------------------------------------------------
$customRange = $d.Content
$customRange.Find.Execute("Custom file format for specifying patterns.", $true, $false, `
                           $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$customPara = $customRange.Paragraphs(1)
$customPara.Range.Delete()
